$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 386.91666
$ws.Range("I12").Value = 294.3
$ws.Range("K12").Value = 294.3
$ws.Range("M12").Value = -124.3
$ws.Range("H70").Value = 3394.45
$ws.Range("I70").Value = 1875
$ws.Range("J70").Value = 3774.3125
$ws.Range("K70").Value = 5625
$ws.Range("L70").Value = 11322.9375
$ws.Range("M70").Value = -5355
$ws.Range("N70").Value = -11862.9375
$ws.Range("H73").Value = 3394.45
$ws.Range("I73").Value = 1875
$ws.Range("J73").Value = 3774.3125
$ws.Range("K73").Value = 5625
$ws.Range("L73").Value = 11322.9375
$ws.Range("M73").Value = -4689
$ws.Range("N73").Value = -13194.9375
$ws.Range("H98").Value = 858.2143
$ws.Range("I98").Value = 931.6667
$ws.Range("J98").Value = 417.5
$ws.Range("K98").Value = 931.6667
$ws.Range("L98").Value = 417.5
$ws.Range("M98").Value = 566.3333
$ws.Range("N98").Value = -3413.5
$ws.Range("H107").Value = 5359.7915
$ws.Range("J107").Value = 8751.375
$ws.Range("L107").Value = 8751.375
$ws.Range("N107").Value = -12591.375
$ws.Range("H122").Value = 858.2143
$ws.Range("I122").Value = 931.6667
$ws.Range("J122").Value = 417.5
$ws.Range("K122").Value = 2795.0001
$ws.Range("L122").Value = 1252.5
$ws.Range("M122").Value = -345.0001000000002
$ws.Range("N122").Value = -6152.5
$ws.Range("H127").Value = 688918.2
$ws.Range("I127").Value = 786942.2
$ws.Range("K127").Value = 2360826.6
$ws.Range("M127").Value = -2355866.6
$ws.Range("H129").Value = 1177521.1
$ws.Range("I129").Value = 667690.7
$ws.Range("K129").Value = 2003072.1
$ws.Range("M129").Value = -1998072.1
$ws.Range("H132").Value = 1343.4762
$ws.Range("I132").Value = 1160.35
$ws.Range("J132").Value = 5006
$ws.Range("K132").Value = 3481.05
$ws.Range("L132").Value = 15018
$ws.Range("M132").Value = -951.0499999999997
$ws.Range("N132").Value = -20078
$ws.Range("H137").Value = 11766114
$ws.Range("I137").Value = 15386011
$ws.Range("K137").Value = 46158033
$ws.Range("M137").Value = -46155483
$ws.Range("H138").Value = 4526.579
$ws.Range("I138").Value = 3731.5
$ws.Range("J138").Value = 4738.6
$ws.Range("K138").Value = 11194.5
$ws.Range("L138").Value = 14215.8
$ws.Range("M138").Value = -6054.5
$ws.Range("N138").Value = -24495.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1879.2858
$ws.Range("I2").Value = 1879.2858
$ws.Range("K2").Value = 1879.2858
$ws.Range("M2").Value = -1766.2858
$ws.Range("H32").Value = 939801.5600000001
$ws.Range("I32").Value = 1203489.9
$ws.Range("J32").Value = 2243.111
$ws.Range("K32").Value = 1203489.9
$ws.Range("L32").Value = 2243.111
$ws.Range("M32").Value = -1203202.9
$ws.Range("N32").Value = -2817.111
$ws.Range("H61").Value = 5404198.5
$ws.Range("I61").Value = 2022419.1
$ws.Range("J61").Value = 27144210
$ws.Range("K61").Value = 2022419.1
$ws.Range("L61").Value = 27144210
$ws.Range("M61").Value = -2022207.1
$ws.Range("N61").Value = -27144634
$ws.Range("H116").Value = 1879.2858
$ws.Range("I116").Value = 1879.2858
$ws.Range("K116").Value = 1879.2858
$ws.Range("M116").Value = 414.7141999999999
$ws.Range("H132").Value = 3263.6584
$ws.Range("I132").Value = 1711.8462
$ws.Range("K132").Value = 5135.5386
$ws.Range("M132").Value = -2605.5386
$ws.Range("H136").Value = 5404198.5
$ws.Range("I136").Value = 2022419.1
$ws.Range("J136").Value = 27144210
$ws.Range("K136").Value = 6067257.300000001
$ws.Range("L136").Value = 81432630
$ws.Range("M136").Value = -6064707.300000001
$ws.Range("N136").Value = -81437730

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1879.2858
$ws.Range("I3").Value = 1879.2858
$ws.Range("K3").Value = 1879.2858
$ws.Range("M3").Value = -1765.2858
$ws.Range("H86").Value = 3790.4614
$ws.Range("I86").Value = 3477.889
$ws.Range("J86").Value = 4493.75
$ws.Range("K86").Value = 3477.889
$ws.Range("L86").Value = 4493.75
$ws.Range("M86").Value = -2354.889
$ws.Range("N86").Value = -6739.75
$ws.Range("H89").Value = 3790.4614
$ws.Range("I89").Value = 3477.889
$ws.Range("J89").Value = 4493.75
$ws.Range("K89").Value = 17389.445
$ws.Range("L89").Value = 22468.75
$ws.Range("M89").Value = -11773.445
$ws.Range("N89").Value = -33700.75
$ws.Range("H99").Value = 15111.889
$ws.Range("I99").Value = 19582.666
$ws.Range("K99").Value = 19582.666
$ws.Range("M99").Value = -18084.666
$ws.Range("H134").Value = 6738354.5
$ws.Range("I134").Value = 5750813.5
$ws.Range("K134").Value = 17252440.5
$ws.Range("M134").Value = -17249905.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 13579.25
$ws.Range("I14").Value = 17606
$ws.Range("K14").Value = 17606
$ws.Range("M14").Value = -17436
$ws.Range("H25").Value = 911
$ws.Range("I25").Value = 911
$ws.Range("K25").Value = 911
$ws.Range("M25").Value = -737
$ws.Range("H31").Value = 532957.9
$ws.Range("I31").Value = 1290657.5
$ws.Range("K31").Value = 1290657.5
$ws.Range("M31").Value = -1290362.5
$ws.Range("H34").Value = 532957.9
$ws.Range("I34").Value = 1290657.5
$ws.Range("K34").Value = 1290657.5
$ws.Range("M34").Value = -1290455.5
$ws.Range("H58").Value = 13413978
$ws.Range("I58").Value = 18522412
$ws.Range("K58").Value = 18522412
$ws.Range("M58").Value = -18522209
$ws.Range("H62").Value = 3093.75
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 3093.75
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 3093.75
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -4341.75
$ws.Range("H65").Value = 3093.75
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 3093.75
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 15468.75
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -21708.75
$ws.Range("H99").Value = 14492.667
$ws.Range("I99").Value = 18980.834
$ws.Range("J99").Value = 5516.3335
$ws.Range("K99").Value = 18980.834
$ws.Range("L99").Value = 5516.3335
$ws.Range("M99").Value = -17482.834
$ws.Range("N99").Value = -8512.333500000001
$ws.Range("H119").Value = 83881
$ws.Range("J119").Value = 83881
$ws.Range("L119").Value = 83881
$ws.Range("N119").Value = -93557
$ws.Range("H126").Value = 14492.667
$ws.Range("I126").Value = 18980.834
$ws.Range("J126").Value = 5516.3335
$ws.Range("K126").Value = 56942.50199999999
$ws.Range("L126").Value = 16549.0005
$ws.Range("M126").Value = -54472.50199999999
$ws.Range("N126").Value = -21489.0005
$ws.Range("H132").Value = 3157
$ws.Range("I132").Value = 2860.25
$ws.Range("K132").Value = 8580.75
$ws.Range("M132").Value = -6050.75
$ws.Range("H136").Value = 13413978
$ws.Range("I136").Value = 18522412
$ws.Range("K136").Value = 55567236
$ws.Range("M136").Value = -55564686

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 3703760.2
$ws.Range("I2").Value = 5555605
$ws.Range("K2").Value = 5555605
$ws.Range("M2").Value = -5555492
$ws.Range("H102").Value = 1636.7894
$ws.Range("I102").Value = 1505.5
$ws.Range("K102").Value = 1505.5
$ws.Range("M102").Value = 116.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H23").Value = 651.75
$ws.Range("I23").Value = 651.75
$ws.Range("K23").Value = 651.75
$ws.Range("M23").Value = -421.75
$ws.Range("H95").Value = 22500
$ws.Range("J95").Value = 22500
$ws.Range("L95").Value = 22500
$ws.Range("N95").Value = -27992

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5147.0835
$ws.Range("I81").Value = 4338.3
$ws.Range("K81").Value = 8676.6
$ws.Range("M81").Value = -7615.6
$ws.Range("H84").Value = 5147.0835
$ws.Range("I84").Value = 4338.3
$ws.Range("K84").Value = 43383
$ws.Range("M84").Value = -38079
$ws.Range("H122").Value = 77664.266
$ws.Range("I122").Value = 2608
$ws.Range("K122").Value = 7824
$ws.Range("M122").Value = -5374
$ws.Range("H126").Value = 2818.6
$ws.Range("I126").Value = 2742.889
$ws.Range("K126").Value = 8228.667000000001
$ws.Range("M126").Value = -5758.667000000001
$ws.Range("H132").Value = 12825075
$ws.Range("J132").Value = 5000
$ws.Range("L132").Value = 15000
$ws.Range("N132").Value = -20060
$ws.Range("H136").Value = 2137556
$ws.Range("I136").Value = 966799.4
$ws.Range("J136").Value = 6950666.5
$ws.Range("K136").Value = 2900398.2
$ws.Range("L136").Value = 20851999.5
$ws.Range("M136").Value = -2897848.2
$ws.Range("N136").Value = -20857099.5
